$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Survey" (sheet6) - lookup lists used by the data validations on the
# "Driver View Test Data" sheet. Done first so the later selection on
# "Driver View Test Data" is the one that sticks as the active sheet/tab.
# ---------------------------------------------------------------------------
$sv = $wb.Worksheets.Item("Survey")

$sv.Cells.Item(1, 5).Value = "Cloud Cover"

$sv.Cells.Item(2, 3).Value = "Calm"
$sv.Cells.Item(2, 5).Value = "LessThan50"

$sv.Cells.Item(3, 2).Value = "Moderate"
$sv.Cells.Item(3, 3).Value = "Light"
$sv.Cells.Item(3, 4).Value = "RapidResponse"
$sv.Cells.Item(3, 5).Value = "GreaterThan50"

$sv.Cells.Item(4, 2).Value = "Strong"
$sv.Cells.Item(4, 3).Value = "Strong"
$sv.Cells.Item(4, 4).Value = "Operator"

$sv.Cells.Item(5, 4).Value = "Manual"

$sv.Cells.Item(6, 4).Value = "Assessment"

$null = $sv.Range("D7").Select()

# ---------------------------------------------------------------------------
# Sheet "Driver View Test Data" (sheet2)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Driver View Test Data")

# Insert a new column at E so the existing Wind/Survey Type/Replay* columns
# slide right one slot (F/G/H/I) keeping their exact widths, then rewrite the
# header row from scratch so the new "Cloud Cover" column lands at F (with
# Wind staying put at E).
$null = $ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666

$ws.Cells.Item(1, 1).Value = "RowID"
$ws.Cells.Item(1, 2).Value = "Survey Tag"
$ws.Cells.Item(1, 3).Value = "Survey Time"
$ws.Cells.Item(1, 4).Value = "Solar Radiation"
$ws.Cells.Item(1, 5).Value = "Wind"
$ws.Cells.Item(1, 6).Value = "Cloud Cover"
$ws.Cells.Item(1, 7).Value = "Survey Type"
$ws.Cells.Item(1, 8).Value = "Replay Script DB3 File"
$ws.Cells.Item(1, 9).Value = "Replay Script Defn File"

# Existing row 2 used to hold the single sample row; it now only keeps the
# RowID and gets a Replay Script Defn File sample value in column I.
$null = $ws.Range("B2:I2").ClearContents()
$ws.Cells.Item(2, 9).Value = "instr_ready.defn"

# Row 3: another defn-file-only sample row.
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 9).Value = "instr_warming.defn"

# Row 4: full sample survey row (previously row 2's data, now reshuffled and
# extended with the new Cloud Cover column).
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "SurveyTag001"
$ws.Cells.Item(4, 3).Value = "Day"
$ws.Cells.Item(4, 4).Value = "Moderate"
$ws.Cells.Item(4, 5).Value = "Light"
$ws.Cells.Item(4, 6).Value = "LessThan50"
$ws.Cells.Item(4, 7).Value = "Standard"
$ws.Cells.Item(4, 8).Value = "Surveyor.db3"
$ws.Cells.Item(4, 9).Value = "replay-db3.defn"

$null = $ws.Range("E17").Select()
